$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 598.3333
$ws.Range("J2").Value = 598.3333
$ws.Range("L2").Value = 598.3333
$ws.Range("N2").Value = -824.3333
$ws.Range("H15").Value = 441.9375
$ws.Range("I15").Value = 441.9375
$ws.Range("K15").Value = 1325.8125
$ws.Range("M15").Value = -1156.8125
$ws.Range("H33").Value = 240.41667
$ws.Range("I33").Value = 93
$ws.Range("J33").Value = 446.8
$ws.Range("K33").Value = 93
$ws.Range("L33").Value = 446.8
$ws.Range("M33").Value = 136
$ws.Range("N33").Value = -904.8
$ws.Range("H51").Value = 7999.6665
$ws.Range("I51").Value = 7999.6665
$ws.Range("K51").Value = 7999.6665
$ws.Range("M51").Value = -7515.6665
$ws.Range("H62").Value = 4889.3
$ws.Range("H65").Value = 4889.3
$ws.Range("H106").Value = 15910.533
$ws.Range("I106").Value = 17719.924
$ws.Range("K106").Value = 17719.924
$ws.Range("M106").Value = -17088.924
$ws.Range("H112").Value = 2303.2222
$ws.Range("I112").Value = 766.3333
$ws.Range("J112").Value = 2610.6
$ws.Range("K112").Value = 2298.9999
$ws.Range("L112").Value = 7831.799999999999
$ws.Range("M112").Value = -1190.9999
$ws.Range("N112").Value = -10047.8
$ws.Range("H116").Value = 5640.8125
$ws.Range("I116").Value = 5378.7
$ws.Range("J116").Value = 6077.6665
$ws.Range("K116").Value = 5378.7
$ws.Range("L116").Value = 6077.6665
$ws.Range("M116").Value = -1936.7
$ws.Range("N116").Value = -12961.6665
$ws.Range("H137").Value = 2319.2942
$ws.Range("I137").Value = 1137
$ws.Range("K137").Value = 3411
$ws.Range("M137").Value = -861
$ws.Range("H138").Value = 3874.7847
$ws.Range("I138").Value = 1271
$ws.Range("J138").Value = 4189.0347
$ws.Range("K138").Value = 3813
$ws.Range("L138").Value = 12567.1041
$ws.Range("M138").Value = 1327
$ws.Range("N138").Value = -22847.1041

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1381.25
$ws.Range("I102").Value = 890
$ws.Range("J102").Value = 2200
$ws.Range("K102").Value = 890
$ws.Range("L102").Value = 2200
$ws.Range("M102").Value = 732
$ws.Range("N102").Value = -5444

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1513.7894
$ws.Range("I20").Value = 729.53845
$ws.Range("K20").Value = 729.53845
$ws.Range("M20").Value = -482.53845
$ws.Range("H105").Value = 4175.0835
$ws.Range("I105").Value = 3313.3157
$ws.Range("K105").Value = 3313.3157
$ws.Range("M105").Value = -1566.3157
$ws.Range("H134").Value = 2582.375
$ws.Range("I134").Value = 2448.842
$ws.Range("K134").Value = 7346.526
$ws.Range("M134").Value = -4811.526

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3628.4167
$ws.Range("I31").Value = 1958.125
$ws.Range("K31").Value = 1958.125
$ws.Range("M31").Value = -1663.125
$ws.Range("H34").Value = 3628.4167
$ws.Range("I34").Value = 1958.125
$ws.Range("K34").Value = 1958.125
$ws.Range("M34").Value = -1756.125
$ws.Range("H58").Value = 3465.611
$ws.Range("I58").Value = 2828.25
$ws.Range("J58").Value = 3647.7144
$ws.Range("K58").Value = 2828.25
$ws.Range("L58").Value = 3647.7144
$ws.Range("M58").Value = -2625.25
$ws.Range("N58").Value = -4053.7144
$ws.Range("H105").Value = 2463.842
$ws.Range("I105").Value = 619.5454999999999
$ws.Range("J105").Value = 4999.75
$ws.Range("K105").Value = 619.5454999999999
$ws.Range("L105").Value = 4999.75
$ws.Range("M105").Value = 1127.4545
$ws.Range("N105").Value = -8493.75
$ws.Range("H107").Value = 1217.1111
$ws.Range("I107").Value = 997.8570999999999
$ws.Range("J107").Value = 1984.5
$ws.Range("K107").Value = 997.8570999999999
$ws.Range("L107").Value = 1984.5
$ws.Range("M107").Value = 922.1429000000001
$ws.Range("N107").Value = -5824.5
$ws.Range("H122").Value = 3906.4546
$ws.Range("I122").Value = 4291.143
$ws.Range("J122").Value = 3233.25
$ws.Range("K122").Value = 12873.429
$ws.Range("L122").Value = 9699.75
$ws.Range("M122").Value = -10423.429
$ws.Range("N122").Value = -14599.75
$ws.Range("H132").Value = 3458.95
$ws.Range("I132").Value = 2182.923
$ws.Range("K132").Value = 6548.768999999999
$ws.Range("M132").Value = -4018.768999999999
$ws.Range("H136").Value = 3465.611
$ws.Range("I136").Value = 2828.25
$ws.Range("J136").Value = 3647.7144
$ws.Range("K136").Value = 8484.75
$ws.Range("L136").Value = 10943.1432
$ws.Range("M136").Value = -5934.75
$ws.Range("N136").Value = -16043.1432

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39289220
$ws.Range("I4").Value = 53763376
$ws.Range("J4").Value = 2226.4285
$ws.Range("K4").Value = 161290128
$ws.Range("L4").Value = 6679.2855
$ws.Range("M4").Value = -161290016
$ws.Range("N4").Value = -6903.2855
$ws.Range("H7").Value = 7692430
$ws.Range("I7").Value = 16666736
$ws.Range("J7").Value = 167.85715
$ws.Range("K7").Value = 50000208
$ws.Range("L7").Value = 503.57145
$ws.Range("M7").Value = -50000096
$ws.Range("N7").Value = -727.5714499999999
$ws.Range("H9").Value = 750
$ws.Range("J9").Value = 750
$ws.Range("L9").Value = 2250
$ws.Range("N9").Value = -2698
$ws.Range("H41").Value = 668.3333
$ws.Range("I41").Value = 655
$ws.Range("J41").Value = 675
$ws.Range("K41").Value = 1965
$ws.Range("L41").Value = 2025
$ws.Range("M41").Value = -1627
$ws.Range("N41").Value = -2701
$ws.Range("H92").Value = 567
$ws.Range("I92").Value = 567
$ws.Range("K92").Value = 1701
$ws.Range("M92").Value = -453
$ws.Range("H100").Value = 2450
$ws.Range("J100").Value = 2450
$ws.Range("L100").Value = 7350
$ws.Range("N100").Value = -8972
$ws.Range("H105").Value = 6325
$ws.Range("J105").Value = 6325
$ws.Range("L105").Value = 18975
$ws.Range("N105").Value = -24217
$ws.Range("H122").Value = 1099.5
$ws.Range("I122").Value = 699
$ws.Range("K122").Value = 6291
$ws.Range("M122").Value = -3841
$ws.Range("H128").Value = 240899.5
$ws.Range("I128").Value = 240899.5
$ws.Range("K128").Value = 722698.5
$ws.Range("M128").Value = -717718.5
$ws.Range("H131").Value = 2407.4707
$ws.Range("J131").Value = 3580.875
$ws.Range("L131").Value = 10742.625
$ws.Range("N131").Value = -20822.625

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1624.625
$ws.Range("I80").Value = 689.6
$ws.Range("J80").Value = 3183
$ws.Range("K80").Value = 689.6
$ws.Range("L80").Value = 3183
$ws.Range("M80").Value = 308.4
$ws.Range("N80").Value = -5179
$ws.Range("H83").Value = 1624.625
$ws.Range("I83").Value = 689.6
$ws.Range("J83").Value = 3183
$ws.Range("K83").Value = 3448
$ws.Range("L83").Value = 15915
$ws.Range("M83").Value = 1544
$ws.Range("N83").Value = -25899
$ws.Range("H102").Value = 1311.7858
$ws.Range("I102").Value = 351.36365
$ws.Range("K102").Value = 351.36365
$ws.Range("M102").Value = 1270.63635
$ws.Range("H132").Value = 2588.3125
$ws.Range("I132").Value = 1710.5
$ws.Range("K132").Value = 5131.5
$ws.Range("M132").Value = -2601.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1441.2941
$ws.Range("I16").Value = 1343.875
$ws.Range("K16").Value = 1343.875
$ws.Range("M16").Value = -1173.875
$ws.Range("H136").Value = 3450.1
$ws.Range("I136").Value = 3628.7144
$ws.Range("K136").Value = 10886.1432
$ws.Range("M136").Value = -8336.143199999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 16000
$ws.Range("I55").Value = 16000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 16000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -15723
$ws.Range("N55").ClearContents()
$ws.Range("H100").Value = 2277.7778
$ws.Range("I100").Value = 2580
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 5160
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -4619
$ws.Range("N100").Value = -4882
$ws.Range("H132").Value = 1322.875
$ws.Range("I132").Value = 1355.6364
$ws.Range("K132").Value = 4066.9092
$ws.Range("M132").Value = -1536.9092
$ws.Range("H136").Value = 3291.2778
$ws.Range("I136").Value = 1421.1538
$ws.Range("J136").Value = 8153.6
$ws.Range("K136").Value = 4263.4614
$ws.Range("L136").Value = 24460.8
$ws.Range("M136").Value = -1713.4614
$ws.Range("N136").Value = -29560.8
